$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegSpec_Example")

# The register bit-field table used plain "hi-lo" bit ranges (e.g. "32-24")
# for the RESERVED fields. Update them to the new multi-range bracket
# notation used by the RTL generator (e.g. "[32:28]").
$ws.Range("A14").Value = "[32:28]"
$ws.Range("A15").Value = "[27:24][23:16]"
$ws.Range("A17").Value = "[14:8][7:1]"

# Bring the sheet to the front and scroll/select the same way the author
# left it (top-left cell A10, active cell A18).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A18").Select()
